$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.618.18"
$ws.Range("E2").Value = "  -1.27%  "
$ws.Range("D3").Value = "1.632.41"
$ws.Range("E3").Value = "  -0.62%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.522"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.83%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.08"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.11%  "
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("E10").Value = "  -0.28%  "
$ws.Range("E11").Value = "  -3.21%  "
$ws.Range("D12").Value = "1.861.96"
$ws.Range("E12").Value = "  -0.76%  "
$ws.Range("D13").Value = "1.631.15"
$ws.Range("E13").Value = "  -0.70%  "
$ws.Range("E14").Value = "  -0.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.559"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.39%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.09"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.62%  "
$ws.Range("D17").Value = "27.588.60"
$ws.Range("E17").Value = "  -1.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "229.39"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.77%  "
$ws.Range("D19").Value = "0.0₃0718"
$ws.Range("E19").Value = "  -0.71%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.56"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.41%  "
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.66"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.62%  "
$ws.Range("E23").Value = "  +1.02%  "
$ws.Range("E24").Value = "  +2.60%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.17"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.93%  "
$ws.Range("E26").Value = "  -1.20%  "
$ws.Range("E27").Value = "  -0.66%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.60"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.56%  "
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("E30").Value = "  -0.88%  "
$ws.Range("E31").Value = "  -0.88%  "
$ws.Range("E32").Value = "  -1.33%  "
$ws.Range("D33").Value = "1.463.77"
$ws.Range("E33").Value = "  -0.47%  "
$ws.Range("E34").Value = "  -0.80%  "
$ws.Range("E35").Value = "  -0.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.32"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.75%  "
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("E38").Value = "  -1.82%  "
$ws.Range("E39").Value = "  -0.67%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.916"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "68.88"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.10%  "
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("E43").Value = "  +0.42%  "
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.38"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.62%  "
$ws.Range("B46").Value = "MXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.21"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.19%  "
$ws.Range("D47").Value = "1.772.19"
$ws.Range("E47").Value = "  -0.83%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.22%  "
$ws.Range("E49").Value = "  +1.10%  "
$ws.Range("E50").Value = "  -0.81%  "
$ws.Range("E51").Value = "  +0.25%  "
